# Add "MaxContinuousRetryNumber" setting to the Constants sheet, and
# tweak the UI selection/active-tab state to match the authored commit
# ("Add abort on MaxContinuousRetryNumber").

$wb = $excel.ActiveWorkbook

$wsConstants = $wb.Worksheets.Item("Constants")
$wsSettings  = $wb.Worksheets.Item("Settings")
$wsIntro     = $wb.Worksheets.Item("Introduction")

# --- Constants sheet: insert a new row above the existing
# "ExScreenshotsFolderPath" row (row 5) for the new retry setting ---
$wsConstants.Rows.Item(5).Insert()

# Write the new setting's values first ...
$wsConstants.Cells.Item(5, 1).Value = "MaxContinuousRetryNumber"
$wsConstants.Cells.Item(5, 2).Value = 0
$wsConstants.Cells.Item(5, 3).Value = "If > 0 will keep a record of consecutive failed exceptions of the Process state. When this number is reached, the application will fail. Must be an integer."

# ... then copy over the formatting of the "MaxInitRetryNumber" row
# above it (row 4), matching Excel's own "format from row above" on
# insert, without disturbing the values just written.
$wsConstants.Range("A4:C4").Copy()
$wsConstants.Range("A5:C5").PasteSpecial(-4122)
$wsConstants.Rows.Item(5).RowHeight = $wsConstants.Rows.Item(4).RowHeight

$wsConstants.Range("B12").Select()

# --- Settings sheet: move the active selection ---
$wsSettings.Range("B38").Select()

# --- Introduction sheet becomes the active tab/sheet ---
$wsIntro.Activate()
$wsIntro.Range("A16").Select()
